$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain TEXT (matching the workbook's original
# inline-string / shared-string cell type) regardless of whether the
# string happens to look numeric (e.g. "0.613", "588.53"). A leading
# apostrophe forces Excel to store it as text instead of auto-converting
# to a number; resetting the Style back to "Normal" afterwards strips the
# transient "quote prefix" cell style so no stray style index is left on
# the cell (matching cells in this sheet, which carry no explicit style).
function Set-TextValue {
    param($cell, [string]$val)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

# --- Row 31 / 32 swap: Filecoin <-> Stacks (with updated price/volume) ---
Set-TextValue $ws.Cells.Item(31,2) "Stacks"
Set-TextValue $ws.Cells.Item(31,3) "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Cells.Item(31,4) "3.14"
Set-TextValue $ws.Cells.Item(31,5) "  -6.29%  "
Set-TextValue $ws.Cells.Item(32,2) "Filecoin"
Set-TextValue $ws.Cells.Item(32,3) "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Cells.Item(32,4) "8.53"
Set-TextValue $ws.Cells.Item(32,5) "  -5.09%  "

# --- Price (D) / Volume(1h) (E) updates for all other rows ---
Set-TextValue $ws.Cells.Item(2,4) "68.960.23"
Set-TextValue $ws.Cells.Item(2,5) "  -2.38%  "
Set-TextValue $ws.Cells.Item(3,4) "3.513.27"
Set-TextValue $ws.Cells.Item(3,5) "  -3.38%  "
Set-TextValue $ws.Cells.Item(4,5) "  -0.02%  "
Set-TextValue $ws.Cells.Item(5,4) "588.53"
Set-TextValue $ws.Cells.Item(5,5) "  +1.24%  "
Set-TextValue $ws.Cells.Item(6,4) "170.88"
Set-TextValue $ws.Cells.Item(6,5) "  -2.70%  "
Set-TextValue $ws.Cells.Item(7,4) "0.613"
Set-TextValue $ws.Cells.Item(7,5) "  +0.52%  "
Set-TextValue $ws.Cells.Item(8,4) "3.505.17"
Set-TextValue $ws.Cells.Item(8,5) "  -3.31%  "
Set-TextValue $ws.Cells.Item(9,5) "  -0.03%  "
Set-TextValue $ws.Cells.Item(10,4) "0.188"
Set-TextValue $ws.Cells.Item(10,5) "  -4.38%  "
Set-TextValue $ws.Cells.Item(11,5) "  +0.08%  "
Set-TextValue $ws.Cells.Item(12,4) "0.581"
Set-TextValue $ws.Cells.Item(12,5) "  -4.05%  "
Set-TextValue $ws.Cells.Item(13,4) "47.19"
Set-TextValue $ws.Cells.Item(13,5) "  -2.50%  "
Set-TextValue $ws.Cells.Item(14,5) "  -3.25%  "
Set-TextValue $ws.Cells.Item(15,4) "4.072.37"
Set-TextValue $ws.Cells.Item(15,5) "  -3.58%  "
Set-TextValue $ws.Cells.Item(16,4) "625.97"
Set-TextValue $ws.Cells.Item(16,5) "  -6.04%  "
Set-TextValue $ws.Cells.Item(17,4) "8.47"
Set-TextValue $ws.Cells.Item(17,5) "  -4.88%  "
Set-TextValue $ws.Cells.Item(18,4) "3.517.43"
Set-TextValue $ws.Cells.Item(18,5) "  -3.01%  "
Set-TextValue $ws.Cells.Item(19,4) "69.041.43"
Set-TextValue $ws.Cells.Item(19,5) "  -2.26%  "
Set-TextValue $ws.Cells.Item(20,4) "0.122"
Set-TextValue $ws.Cells.Item(20,5) "  +0.26%  "
Set-TextValue $ws.Cells.Item(21,4) "17.39"
Set-TextValue $ws.Cells.Item(21,5) "  -2.16%  "
Set-TextValue $ws.Cells.Item(22,4) "11.13"
Set-TextValue $ws.Cells.Item(22,5) "  -2.37%  "
Set-TextValue $ws.Cells.Item(23,4) "0.886"
Set-TextValue $ws.Cells.Item(23,5) "  -5.64%  "
Set-TextValue $ws.Cells.Item(24,4) "15.92"
Set-TextValue $ws.Cells.Item(24,5) "  -7.01%  "
Set-TextValue $ws.Cells.Item(25,4) "96.93"
Set-TextValue $ws.Cells.Item(25,5) "  -2.78%  "
Set-TextValue $ws.Cells.Item(26,5) "  -2.43%  "
Set-TextValue $ws.Cells.Item(27,5) "  -0.03%  "
Set-TextValue $ws.Cells.Item(28,5) "  -5.64%  "
Set-TextValue $ws.Cells.Item(29,4) "9.25"
Set-TextValue $ws.Cells.Item(29,5) "  -7.01%  "
Set-TextValue $ws.Cells.Item(30,4) "32.74"
Set-TextValue $ws.Cells.Item(30,5) "  -5.63%  "
Set-TextValue $ws.Cells.Item(33,4) "1.32"
Set-TextValue $ws.Cells.Item(33,5) "  -5.05%  "
Set-TextValue $ws.Cells.Item(34,5) "  -6.38%  "
Set-TextValue $ws.Cells.Item(35,4) "639.30"
Set-TextValue $ws.Cells.Item(35,5) "  +8.83%  "
Set-TextValue $ws.Cells.Item(36,5) "  -2.84%  "
Set-TextValue $ws.Cells.Item(37,4) "3.48"
Set-TextValue $ws.Cells.Item(37,5) "  -13.13%  "
Set-TextValue $ws.Cells.Item(38,4) "0.102"
Set-TextValue $ws.Cells.Item(38,5) "  -4.25%  "
Set-TextValue $ws.Cells.Item(39,4) "57.24"
Set-TextValue $ws.Cells.Item(39,5) "  -1.71%  "
Set-TextValue $ws.Cells.Item(40,5) "  +0.20%  "
Set-TextValue $ws.Cells.Item(41,4) "0.0453"
Set-TextValue $ws.Cells.Item(41,5) "  -0.49%  "
Set-TextValue $ws.Cells.Item(42,5) "  -4.09%  "
Set-TextValue $ws.Cells.Item(43,4) "3.379.63"
Set-TextValue $ws.Cells.Item(43,5) "  -5.37%  "
Set-TextValue $ws.Cells.Item(44,4) "0.327"
Set-TextValue $ws.Cells.Item(44,5) "  -4.82%  "
Set-TextValue $ws.Cells.Item(45,4) "32.79"
Set-TextValue $ws.Cells.Item(45,5) "  -4.90%  "
Set-TextValue $ws.Cells.Item(46,4) "0.0₃0697"
Set-TextValue $ws.Cells.Item(46,5) "  -5.68%  "
Set-TextValue $ws.Cells.Item(47,4) "2.54"
Set-TextValue $ws.Cells.Item(47,5) "  -5.62%  "
Set-TextValue $ws.Cells.Item(48,5) "  -6.11%  "
Set-TextValue $ws.Cells.Item(49,5) "  -2.20%  "
Set-TextValue $ws.Cells.Item(50,5) "  -2.26%  "
Set-TextValue $ws.Cells.Item(51,5) "  +14.63%  "
